$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header cell formatting (bold, border, alignment) from E1 onto F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$times = @(
    "2021-10-05 10:51:05.844865",
    "2021-10-05 10:51:05.844876",
    "2021-10-05 10:51:05.844879",
    "2021-10-05 10:51:05.844882",
    "2021-10-05 10:51:05.844885",
    "2021-10-05 10:51:05.844888",
    "2021-10-05 10:51:05.844890",
    "2021-10-05 10:51:05.844893",
    "2021-10-05 10:51:05.844896",
    "2021-10-05 10:51:05.844898",
    "2021-10-05 10:51:05.844901",
    "2021-10-05 10:51:05.844903",
    "2021-10-05 10:51:05.844906",
    "2021-10-05 10:51:05.844909",
    "2021-10-05 10:51:05.844911",
    "2021-10-05 10:51:05.844914",
    "2021-10-05 10:51:05.844917",
    "2021-10-05 10:51:05.844919",
    "2021-10-05 10:51:05.844922",
    "2021-10-05 10:51:05.844925",
    "2021-10-05 10:51:05.844927",
    "2021-10-05 10:51:05.844930",
    "2021-10-05 10:51:05.844932",
    "2021-10-05 10:51:05.844935",
    "2021-10-05 10:51:05.844938",
    "2021-10-05 10:51:05.844941",
    "2021-10-05 10:51:05.844943"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
